# Updates cryptos list values (price + volume/1h) to reflect the latest
# scrape, and fixes the row ordering for three coin pairs whose rank
# swapped (Litecoin/PEPE, Cosmos/Arweave).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into the "Price" column while keeping it stored as
# plain text (the sheet uses text cells for prices, some of which look like
# numbers e.g. "600.41" or have multiple dots e.g. "69.295.31"). Forcing a
# text number format prevents Excel from silently converting the string to
# a floating point number (which would lose formatting / precision), and
# resetting the style back to "Normal" afterwards keeps the cell's style
# identical to its original (unstyled) state.
function Set-PriceText($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-PriceText $ws.Range("D2") "69.295.31"
$ws.Range("E2").Value = "  +2.90%  "

# Row 3 - Ethereum
Set-PriceText $ws.Range("D3") "3.806.95"
$ws.Range("E3").Value = "  +1.44%  "

# Row 4 - TetherUSD
Set-PriceText $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
Set-PriceText $ws.Range("D5") "600.41"
$ws.Range("E5").Value = "  +0.97%  "

# Row 6 - Solana
Set-PriceText $ws.Range("D6") "171.04"
$ws.Range("E6").Value = "  +0.60%  "

# Row 7 - LidoStakedEther
Set-PriceText $ws.Range("D7") "3.804.15"
$ws.Range("E7").Value = "  +1.44%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  -0.02%  "

# Row 10 - Dogecoin
Set-PriceText $ws.Range("D10") "0.163"
$ws.Range("E10").Value = "  -1.36%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +0.68%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  -0.21%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  -4.55%  "

# Row 14 - Avalanche
Set-PriceText $ws.Range("D14") "36.92"
$ws.Range("E14").Value = "  +0.73%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-PriceText $ws.Range("D15") "4.446.73"
$ws.Range("E15").Value = "  +1.45%  "

# Row 16 - WrappedEther
Set-PriceText $ws.Range("D16") "3.825.93"
$ws.Range("E16").Value = "  +1.86%  "

# Row 17 - WrappedBTC
Set-PriceText $ws.Range("D17") "69.278.39"
$ws.Range("E17").Value = "  +2.80%  "

# Row 18 - Chainlink
Set-PriceText $ws.Range("D18") "18.26"
$ws.Range("E18").Value = "  -2.11%  "

# Row 19 - Polkadot
Set-PriceText $ws.Range("D19") "7.09"
$ws.Range("E19").Value = "  -1.62%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -0.18%  "

# Row 21 - Uniswap
Set-PriceText $ws.Range("D21") "11.12"
$ws.Range("E21").Value = "  +5.84%  "

# Row 22 - BitcoinCash
Set-PriceText $ws.Range("D22") "472.13"
$ws.Range("E22").Value = "  +0.71%  "

# Row 23 - Polygon
$ws.Range("E23").Value = "  -1.46%  "

# Row 24 - was Litecoin, becomes PEPE
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-PriceText $ws.Range("D24") "0.0000149"
$ws.Range("E24").Value = "  +1.62%  "

# Row 25 - was PEPE, becomes Litecoin
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-PriceText $ws.Range("D25") "84.91"
$ws.Range("E25").Value = "  +1.35%  "

# Row 26 - Fetch.AI
$ws.Range("E26").Value = "  +1.27%  "

# Row 27 - InternetComputer(DFINITY)
Set-PriceText $ws.Range("D27") "12.23"
$ws.Range("E27").Value = "  +0.51%  "

# Row 28 - RenderToken
Set-PriceText $ws.Range("D28") "10.30"
$ws.Range("E28").Value = "  -0.71%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.14%  "

# Row 30 - WrappedeETH
Set-PriceText $ws.Range("D30") "3.955.94"
$ws.Range("E30").Value = "  +1.28%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -2.41%  "

# Row 32 - NEARProtocol
Set-PriceText $ws.Range("D32") "7.48"
$ws.Range("E32").Value = "  -2.75%  "

# Row 33 - ImmutableX
$ws.Range("E33").Value = "  +0.55%  "

# Row 34 - EthereumClassic
Set-PriceText $ws.Range("D34") "30.37"
$ws.Range("E34").Value = "  -0.27%  "

# Row 35 - Aptos
Set-PriceText $ws.Range("D35") "9.43"
$ws.Range("E35").Value = "  +3.26%  "

# Row 37 - RenzoRestakedETH
Set-PriceText $ws.Range("D37") "3.761.80"
$ws.Range("E37").Value = "  +1.19%  "

# Row 38 - Hedera
$ws.Range("E38").Value = "  -2.05%  "

# Row 39 - dogwifhat
Set-PriceText $ws.Range("D39") "3.53"
$ws.Range("E39").Value = "  -7.51%  "

# Row 40 - Kaspa
$ws.Range("E40").Value = "  +1.48%  "

# Row 41 - Mantle
$ws.Range("E41").Value = "  +1.79%  "

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +0.15%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.05%  "

# Row 44 - TheGraph
Set-PriceText $ws.Range("D44") "0.311"
$ws.Range("E44").Value = "  -0.32%  "

# Row 46 - Stacks
$ws.Range("E46").Value = "  +1.48%  "

# Row 47 - was Arweave, becomes Cosmos
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-PriceText $ws.Range("D47") "8.67"
$ws.Range("E47").Value = "  -0.62%  "

# Row 48 - was Cosmos, becomes Arweave
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-PriceText $ws.Range("D48") "43.52"
$ws.Range("E48").Value = "  +10.82%  "

# Row 49 - OKB
Set-PriceText $ws.Range("D49") "46.13"
$ws.Range("E49").Value = "  +0.55%  "

# Row 50 - Bittensor
Set-PriceText $ws.Range("D50") "403.28"
$ws.Range("E50").Value = "  +1.03%  "

# Row 51 - Monero
Set-PriceText $ws.Range("D51") "144.70"
$ws.Range("E51").Value = "  +3.10%  "
